$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) ProcessPayrollForJan16MonthSAPP: bump the report month from 2018 to 2019
# ---------------------------------------------------------------------------
$wsJan = $wb.Worksheets.Item("ProcessPayrollForJan16MonthSAPP")
$wsJan.Range("E2").Value = "January-2019"

# ---------------------------------------------------------------------------
# 2) CreateLeaveRequest: add three new data columns (employeeTaxable,
#    employeeNiable, includeInHolidayEarnings) ahead of the Runmode/Pass/Error
#    trio, each populated with "1" on the data row.
# ---------------------------------------------------------------------------
$wsLeave = $wb.Worksheets.Item("CreateLeaveRequest")

# Shift the existing Runmode/Pass/Error columns (J:L) three slots to the right
# (to M:O) by inserting three blank columns in front of them.
$wsLeave.Columns("J:L").Insert()

# Restore the explicit widths lost on the columns pushed right by the insert.
$wsLeave.Columns("A").ColumnWidth = 37.28515625
$wsLeave.Columns("B").ColumnWidth = 11.7109375
$wsLeave.Columns("C").ColumnWidth = 22
$wsLeave.Columns("D").ColumnWidth = 14.28515625
$wsLeave.Columns("E").ColumnWidth = 14.7109375
$wsLeave.Columns("F").ColumnWidth = 12.7109375
$wsLeave.Columns("G").ColumnWidth = 12.5703125
$wsLeave.Columns("H").ColumnWidth = 38.7109375
$wsLeave.Columns("I").ColumnWidth = 19.7109375
$wsLeave.Columns("J").ColumnWidth = 18
$wsLeave.Columns("K").ColumnWidth = 24.5703125
$wsLeave.Columns("L").ColumnWidth = 24.5703125

# Header row for the new columns.
$wsLeave.Range("J1").Value = "employeeTaxable"
$wsLeave.Range("K1").Value = "employeeNiable"
$wsLeave.Range("L1").Value = "includeInHolidayEarnings"

# Match the header formatting used by the other headers on the row.
$wsLeave.Range("I1").Copy()
$wsLeave.Range("J1:L1").PasteSpecial(-4122)

# Data row values for the new columns.
$wsLeave.Range("J2").Value = "1"
$wsLeave.Range("K2").Value = "1"
$wsLeave.Range("L2").Value = "1"

# Pull the matching cell format from an equivalent report sheet that already
# carries the "1" marker-column style.
$wsReport = $wb.Worksheets.Item("AverageWeeklyEarningsTestReport")
$wsReport.Range("J2").Copy()
$wsLeave.Range("J2:L2").PasteSpecial(-4122)

$wsReport.Range("L2").Copy()
$wsLeave.Range("M2").PasteSpecial(-4122)

$wsReport.Range("M2").Copy()
$wsLeave.Range("N2:O2").PasteSpecial(-4122)
